$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.061.00'
$ws.Range("E2").Value = '  +5.26%  '
$ws.Range("D3").Value = '1.879.29'
$ws.Range("E3").Value = '  +4.00%  '
$ws.Range("E4").Value = '  +0.05%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '280.99'
$c.NumberFormat = "General"
$ws.Range("E5").Value = '  +2.20%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.NumberFormat = "General"
$ws.Range("E6").Value = '  +0.03%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5291'
$c.NumberFormat = "General"
$ws.Range("E7").Value = '  +4.72%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3520'
$c.NumberFormat = "General"
$ws.Range("E8").Value = '  +0.20%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '45.43'
$c.NumberFormat = "General"
$ws.Range("E9").Value = '  +2.53%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.07033'
$c.NumberFormat = "General"
$ws.Range("E10").Value = '  +6.07%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '20.33'
$c.NumberFormat = "General"
$ws.Range("E11").Value = '  +2.09%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.8151'
$c.NumberFormat = "General"
$ws.Range("E12").Value = '  -2.10%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.07807'
$c.NumberFormat = "General"
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("D14").Value = '1.883.30'
$ws.Range("E14").Value = '  +4.27%  '
$ws.Range("E15").Value = '  +2.92%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '90.52'
$c.NumberFormat = "General"
$ws.Range("E16").Value = '  +3.66%  '
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("E18").Value = '  +4.98%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000008203'
$c.NumberFormat = "General"
$ws.Range("E19").Value = '  +2.85%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '1.0000'
$c.NumberFormat = "General"
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").Value = '27.099.25'
$ws.Range("E21").Value = '  +5.18%  '
$ws.Range("D22").Value = '2.118.79'
$ws.Range("E22").Value = '  +4.30%  '
$ws.Range("E23").Value = '  +1.09%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '10.15'
$c.NumberFormat = "General"
$ws.Range("E24").Value = '  +2.02%  '
$ws.Range("E25").Value = '  +2.84%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.386'
$c.NumberFormat = "General"
$ws.Range("E26").Value = '  +12.45%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '146.38'
$c.NumberFormat = "General"
$ws.Range("E27").Value = '  +3.20%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '17.59'
$c.NumberFormat = "General"
$ws.Range("E28").Value = '  +4.04%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.675'
$c.NumberFormat = "General"
$ws.Range("E29").Value = '  +1.50%  '
$ws.Range("E30").Value = '  +3.82%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.392'
$c.NumberFormat = "General"
$ws.Range("E31").Value = '  +1.57%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.373'
$c.NumberFormat = "General"
$ws.Range("E32").Value = '  +4.38%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.08906'
$c.NumberFormat = "General"
$ws.Range("E33").Value = '  +1.55%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.04904'
$c.NumberFormat = "General"
$ws.Range("E34").Value = '  +2.41%  '
$ws.Range("E35").Value = '  +3.74%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.7449'
$c.NumberFormat = "General"
$ws.Range("E36").Value = '  +3.09%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.899'
$c.NumberFormat = "General"
$ws.Range("E37").Value = '  +1.00%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.309'
$c.NumberFormat = "General"
$ws.Range("E38").Value = '  +9.14%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.411'
$c.NumberFormat = "General"
$ws.Range("E39").Value = '  +5.89%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.5319'
$c.NumberFormat = "General"
$ws.Range("E40").Value = '  +2.92%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.01883'
$c.NumberFormat = "General"
$ws.Range("E41").Value = '  +1.56%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.9758'
$c.NumberFormat = "General"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '117.14'
$c.NumberFormat = "General"
$ws.Range("E43").Value = '  +4.06%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '6.315'
$c.NumberFormat = "General"
$ws.Range("E44").Value = '  +2.71%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.204'
$c.NumberFormat = "General"
$ws.Range("E45").Value = '  +2.62%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.9992'
$c.NumberFormat = "General"
$ws.Range("E46").Value = '  +0.03%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.4600'
$c.NumberFormat = "General"
$ws.Range("E47").Value = '  +0.99%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.1369'
$c.NumberFormat = "General"
$ws.Range("E48").Value = '  -0.44%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '9.429'
$c.NumberFormat = "General"
$ws.Range("E49").Value = '  +1.98%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '36.70'
$c.NumberFormat = "General"
$ws.Range("E50").Value = '  +1.94%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.528'
$c.NumberFormat = "General"
$ws.Range("E51").Value = '  +2.58%  '
